$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures at the top of the statement ---
# "VALOR MORA" total (E11): 170820 -> 341640
$ws.Range("E11").Value = 341640
# "Cant. Periodos" (F13): 1 -> 2
$ws.Range("F13").Value = 2

# --- Add a second block of worker rows (period 2508) below the existing
#     period-2507 block, i.e. insert 3 new table rows right after row 18 ---
$ws.Rows("19:21").Insert()

# Row 21 becomes the new "last row" of the table, so give it the bottom-
# border style that row 18 (the old last row) still carries at this point.
$ws.Range("B18:J18").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# Rows 19 and 20 are normal (non-last) rows - copy the plain row style.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)

# Row 18 is no longer the last row of the table, so switch it to the
# normal row style too.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

$ws.CutCopyMode = 0

# Row 19: ANA KARINA PATERNINA WILCHES, period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1235044464"
$ws.Range("D19").Value = "ANA KARINA PATERNINA WILCHES"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20: MARGARETH AVILA CASTAÑO, period 2508
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047432815"
$ws.Range("D20").Value = "MARGARETH AVILA CASTAÑO"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# Row 21: OSCAR JAVIER VASQUEZ REYES, period 2508
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143407895"
$ws.Range("D21").Value = "OSCAR JAVIER VASQUEZ REYES"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500
